# Corrige la ortografía en el documento.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "Rmarkdown" "RStudio"
Replace-Text "Seleccionar version de control" "Seleccionar versión de control"
Replace-Text "pegar el directorio clonado de github y dar click en el boton señalado." "pegar el directorio clonado de github y dar click en el botón señalado."
Replace-Text "A continuacion en la ventana de RStudio darclick en la parte señalada en la siguiente imagen." "A continuación en la ventana de RStudio dar click en la parte señalada en la siguiente imagen."
Replace-Text "En el boton señalado escoger la opcion" "En el botón señalado escoger la opción"
Replace-Text "En la ventana resultante se copian las tras lineas de codigo siguientes:" "En la ventana resultante se copian las tres líneas de código siguientes:"
Replace-Text "Al hacer push por primera vez se debe registar el usuario y la contraseña de la cuenta de github." "Al hacer push por primera vez se debe registrar el usuario y la contraseña de la cuenta de github."
Replace-Text "Cuando Git no esta instalado" "Cuando Git no está instalado"
